$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb2"
$ws.Range("C2").Value = "Ephb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 28.95628266666667
$ws.Range("H2").Value = 86.868848
$ws.Range("I2").Value = 0.5491054194301004
$ws.Range("J2").Value = 0.5491054194301005
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.143611
$ws.Range("N2").Value = 0.430833
$ws.Range("O2").Value = 0.006856337892517759
$ws.Range("P2").Value = 0.006856337892517758
$ws.Range("Q2").Value = 4.158440710042667
$ws.Range("R2").Value = 37.425966390384
$ws.Range("S2").Value = 0.003764852294225455
$ws.Range("T2").Value = 0.003764852294225455
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb2"
$ws.Range("C3").Value = "Ephb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 28.95628266666667
$ws.Range("H3").Value = 86.868848
$ws.Range("I3").Value = 0.5491054194301004
$ws.Range("J3").Value = 0.5491054194301005
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 17.192962
$ws.Range("N3").Value = 51.578886
$ws.Range("O3").Value = 0.8208337581746376
$ws.Range("P3").Value = 0.8208337581746377
$ws.Range("Q3").Value = 497.8442675492586
$ws.Range("R3").Value = 4480.598407943327
$ws.Range("S3").Value = 0.45072426506487
$ws.Range("T3").Value = 0.4507242650648701
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efnb2"
$ws.Range("C4").Value = "Ephb3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 28.95628266666667
$ws.Range("H4").Value = 86.868848
$ws.Range("I4").Value = 0.5491054194301004
$ws.Range("J4").Value = 0.5491054194301005
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.609156666666667
$ws.Range("N4").Value = 10.82747
$ws.Range("O4").Value = 0.1723099039328446
$ws.Range("P4").Value = 0.1723099039328446
$ws.Range("Q4").Value = 104.5077606282844
$ws.Range("R4").Value = 940.56984565456
$ws.Range("S4").Value = 0.09461630207100495
$ws.Range("T4").Value = 0.09461630207100497
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efnb2"
$ws.Range("C5").Value = "Ephb3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 12.691493
$ws.Range("H5").Value = 38.074479
$ws.Range("I5").Value = 0.2406720388519202
$ws.Range("J5").Value = 0.2406720388519202
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.143611
$ws.Range("N5").Value = 0.430833
$ws.Range("O5").Value = 0.006856337892517759
$ws.Range("P5").Value = 0.006856337892517758
$ws.Range("Q5").Value = 1.822638001223
$ws.Range("R5").Value = 16.403742011007
$ws.Range("S5").Value = 0.001650128819649927
$ws.Range("T5").Value = 0.001650128819649926
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efnb2"
$ws.Range("C6").Value = "Ephb3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 12.691493
$ws.Range("H6").Value = 38.074479
$ws.Range("I6").Value = 0.2406720388519202
$ws.Range("J6").Value = 0.2406720388519202
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 17.192962
$ws.Range("N6").Value = 51.578886
$ws.Range("O6").Value = 0.8208337581746376
$ws.Range("P6").Value = 0.8208337581746377
$ws.Range("Q6").Value = 218.204356872266
$ws.Range("R6").Value = 1963.839211850394
$ws.Range("S6").Value = 0.197551734138374
$ws.Range("T6").Value = 0.197551734138374
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efnb2"
$ws.Range("C7").Value = "Ephb3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 12.691493
$ws.Range("H7").Value = 38.074479
$ws.Range("I7").Value = 0.2406720388519202
$ws.Range("J7").Value = 0.2406720388519202
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.609156666666667
$ws.Range("N7").Value = 10.82747
$ws.Range("O7").Value = 0.1723099039328446
$ws.Range("P7").Value = 0.1723099039328446
$ws.Range("Q7").Value = 45.80558657090333
$ws.Range("R7").Value = 412.2502791381299
$ws.Range("S7").Value = 0.04147017589389621
$ws.Range("T7").Value = 0.04147017589389621
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Efnb2"
$ws.Range("C8").Value = "Ephb3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.4888703333333334
$ws.Range("H8").Value = 1.466611
$ws.Range("I8").Value = 0.009270573592685367
$ws.Range("J8").Value = 0.009270573592685367
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.143611
$ws.Range("N8").Value = 0.430833
$ws.Range("O8").Value = 0.006856337892517759
$ws.Range("P8").Value = 0.006856337892517758
$ws.Range("Q8").Value = 0.07020715744033335
$ws.Range("R8").Value = 0.6318644169630001
$ws.Range("S8").Value = 0.00006356218500890318
$ws.Range("T8").Value = 0.00006356218500890318
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Efnb2"
$ws.Range("C9").Value = "Ephb3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.4888703333333334
$ws.Range("H9").Value = 1.466611
$ws.Range("I9").Value = 0.009270573592685367
$ws.Range("J9").Value = 0.009270573592685367
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 17.192962
$ws.Range("N9").Value = 51.578886
$ws.Range("O9").Value = 0.8208337581746376
$ws.Range("P9").Value = 0.8208337581746377
$ws.Range("Q9").Value = 8.405129063927333
$ws.Range("R9").Value = 75.646161575346
$ws.Range("S9").Value = 0.007609599762518481
$ws.Range("T9").Value = 0.007609599762518482
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Efnb2"
$ws.Range("C10").Value = "Ephb3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.4888703333333334
$ws.Range("H10").Value = 1.466611
$ws.Range("I10").Value = 0.009270573592685367
$ws.Range("J10").Value = 0.009270573592685367
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.609156666666667
$ws.Range("N10").Value = 10.82747
$ws.Range("O10").Value = 0.1723099039328446
$ws.Range("P10").Value = 0.1723099039328446
$ws.Range("Q10").Value = 1.764409622685556
$ws.Range("R10").Value = 15.87968660417
$ws.Range("S10").Value = 0.001597411645157982
$ws.Range("T10").Value = 0.001597411645157982
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Efnb2"
$ws.Range("C11").Value = "Ephb3"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 10.59691233333333
$ws.Range("H11").Value = 31.790737
$ws.Range("I11").Value = 0.200951968125294
$ws.Range("J11").Value = 0.200951968125294
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.143611
$ws.Range("N11").Value = 0.430833
$ws.Range("O11").Value = 0.006856337892517759
$ws.Range("P11").Value = 0.006856337892517758
$ws.Range("Q11").Value = 1.521833177102334
$ws.Range("R11").Value = 13.696498593921
$ws.Range("S11").Value = 0.001377794593633474
$ws.Range("T11").Value = 0.001377794593633474
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Efnb2"
$ws.Range("C12").Value = "Ephb3"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 10.59691233333333
$ws.Range("H12").Value = 31.790737
$ws.Range("I12").Value = 0.200951968125294
$ws.Range("J12").Value = 0.200951968125294
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 17.192962
$ws.Range("N12").Value = 51.578886
$ws.Range("O12").Value = 0.8208337581746376
$ws.Range("P12").Value = 0.8208337581746377
$ws.Range("Q12").Value = 182.1923110643313
$ws.Range("R12").Value = 1639.730799578982
$ws.Range("S12").Value = 0.1649481592088751
$ws.Range("T12").Value = 0.1649481592088751
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Efnb2"
$ws.Range("C13").Value = "Ephb3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 10.59691233333333
$ws.Range("H13").Value = 31.790737
$ws.Range("I13").Value = 0.200951968125294
$ws.Range("J13").Value = 0.200951968125294
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 3.609156666666667
$ws.Range("N13").Value = 10.82747
$ws.Range("O13").Value = 0.1723099039328446
$ws.Range("P13").Value = 0.1723099039328446
$ws.Range("Q13").Value = 38.24591679393222
$ws.Range("R13").Value = 344.21325114539
$ws.Range("S13").Value = 0.03462601432278547
$ws.Range("T13").Value = 0.03462601432278547
